# Add a new patient record (row 5) to the Patient_List worksheet.
#
# Columns: A=Patient ID, B=Name, C=Date of Birth, D=Gender,
#          E=Blood Type, F=Phone Number, G=Contact Information

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "P1004"
$ws.Range("B5").Value = "Hi"

# "1232-10-23" looks like a valid calendar date, so a plain assignment
# would get auto-converted into a date serial number by the engine's
# input parser. Route it through a text formula and then paste back
# just the computed value, which keeps it as a genuine shared-string
# cell (matching the rest of the row's formatting/style) instead of
# re-triggering the date auto-detection.
$ws.Range("C5").Formula = '="1232-10-23"'
$ws.Range("C5").Copy()
$ws.Range("C5").PasteSpecial(-4163)

$ws.Range("D5").Value = "Female"
$ws.Range("E5").Value = "O-"
$ws.Range("F5").Value = 98736151
$ws.Range("G5").Value = "1a1"
